# Auto-generated edit script: update market-price derived columns (H-N)
# on the leve-profit tracker sheets, matching the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 321.42856
$ws.Cells.Item(9, 9).Value = 87.5
$ws.Cells.Item(9, 10).Value = 633.3333
$ws.Cells.Item(9, 11).Value = 87.5
$ws.Cells.Item(9, 12).Value = 633.3333
$ws.Cells.Item(9, 13).Value = 81.5
$ws.Cells.Item(9, 14).Value = -971.3333
$ws.Cells.Item(62, 8).Value = 3945
$ws.Cells.Item(62, 9).Value = 2622.8
$ws.Cells.Item(62, 11).Value = 2622.8
$ws.Cells.Item(62, 13).Value = -1998.8
$ws.Cells.Item(65, 8).Value = 3945
$ws.Cells.Item(65, 9).Value = 2622.8
$ws.Cells.Item(65, 11).Value = 13114
$ws.Cells.Item(65, 13).Value = -9994
$ws.Cells.Item(88, 8).Value = 1166.6666
$ws.Cells.Item(88, 9).Value = 1000
$ws.Cells.Item(88, 10).Value = 1500
$ws.Cells.Item(88, 11).Value = 1000
$ws.Cells.Item(88, 12).Value = 1500
$ws.Cells.Item(88, 13).Value = -594
$ws.Cells.Item(88, 14).Value = -2312
$ws.Cells.Item(91, 8).Value = 1166.6666
$ws.Cells.Item(91, 9).Value = 1000
$ws.Cells.Item(91, 10).Value = 1500
$ws.Cells.Item(91, 11).Value = 1000
$ws.Cells.Item(91, 12).Value = 1500
$ws.Cells.Item(91, 13).Value = 404
$ws.Cells.Item(91, 14).Value = -4308
$ws.Cells.Item(99, 8).Value = 235.33333
$ws.Cells.Item(99, 9).Value = 193.6
$ws.Cells.Item(99, 10).Value = 287.5
$ws.Cells.Item(99, 11).Value = 580.8
$ws.Cells.Item(99, 12).Value = 862.5
$ws.Cells.Item(99, 13).Value = 917.2
$ws.Cells.Item(99, 14).Value = -3858.5
$ws.Cells.Item(112, 8).Value = 1059.3469
$ws.Cells.Item(112, 10).Value = 1070.3829
$ws.Cells.Item(112, 12).Value = 3211.1487
$ws.Cells.Item(112, 14).Value = -5427.1487
$ws.Cells.Item(129, 8).Value = 239051.27
$ws.Cells.Item(129, 9).Value = 247.8
$ws.Cells.Item(129, 10).Value = 313677.34
$ws.Cells.Item(129, 11).Value = 743.4000000000001
$ws.Cells.Item(129, 12).Value = 941032.02
$ws.Cells.Item(129, 13).Value = 4256.6
$ws.Cells.Item(129, 14).Value = -951032.02
$ws.Cells.Item(131, 8).Value = 2196.2144
$ws.Cells.Item(131, 9).Value = 1642.4286
$ws.Cells.Item(131, 11).Value = 4927.2858
$ws.Cells.Item(131, 13).Value = 112.7142000000003
$ws.Cells.Item(132, 8).Value = 3409.6428
$ws.Cells.Item(132, 9).Value = 3480.3704
$ws.Cells.Item(132, 11).Value = 10441.1112
$ws.Cells.Item(132, 13).Value = -7911.111199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 727.44446
$ws.Cells.Item(2, 9).Value = 697.7308
$ws.Cells.Item(2, 11).Value = 697.7308
$ws.Cells.Item(2, 13).Value = -584.7308
$ws.Cells.Item(32, 8).Value = 7977.7456
$ws.Cells.Item(32, 9).Value = 6601.431
$ws.Cells.Item(32, 11).Value = 6601.431
$ws.Cells.Item(32, 13).Value = -6314.431
$ws.Cells.Item(45, 8).Value = 4155.2666
$ws.Cells.Item(45, 9).Value = 4100.1113
$ws.Cells.Item(45, 10).Value = 4238
$ws.Cells.Item(45, 11).Value = 4100.1113
$ws.Cells.Item(45, 12).Value = 4238
$ws.Cells.Item(45, 13).Value = -3723.1113
$ws.Cells.Item(45, 14).Value = -4992
$ws.Cells.Item(74, 8).Value = 30304774
$ws.Cells.Item(74, 9).Value = 52632230
$ws.Cells.Item(74, 11).Value = 52632230
$ws.Cells.Item(74, 13).Value = -52631356
$ws.Cells.Item(77, 8).Value = 30304774
$ws.Cells.Item(77, 9).Value = 52632230
$ws.Cells.Item(77, 11).Value = 263161150
$ws.Cells.Item(77, 13).Value = -263156782
$ws.Cells.Item(98, 8).Value = 30000
$ws.Cells.Item(98, 10).Value = 30000
$ws.Cells.Item(98, 12).Value = 30000
$ws.Cells.Item(98, 14).Value = -35990
$ws.Cells.Item(116, 8).Value = 727.44446
$ws.Cells.Item(116, 9).Value = 697.7308
$ws.Cells.Item(116, 11).Value = 697.7308
$ws.Cells.Item(116, 13).Value = 1596.2692

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 727.44446
$ws.Cells.Item(3, 9).Value = 697.7308
$ws.Cells.Item(3, 11).Value = 697.7308
$ws.Cells.Item(3, 13).Value = -583.7308
$ws.Cells.Item(99, 8).Value = 1727.5
$ws.Cells.Item(99, 9).Value = 1633
$ws.Cells.Item(99, 11).Value = 1633
$ws.Cells.Item(99, 13).Value = -135
$ws.Cells.Item(112, 8).Value = 43469.5
$ws.Cells.Item(112, 10).Value = 43469.5
$ws.Cells.Item(112, 12).Value = 43469.5
$ws.Cells.Item(112, 14).Value = -46423.5
$ws.Cells.Item(134, 8).Value = 3758.634
$ws.Cells.Item(134, 9).Value = 3704.4
$ws.Cells.Item(134, 10).Value = 4075
$ws.Cells.Item(134, 11).Value = 11113.2
$ws.Cells.Item(134, 12).Value = 12225
$ws.Cells.Item(134, 13).Value = -8578.200000000001
$ws.Cells.Item(134, 14).Value = -17295

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3008.0588
$ws.Cells.Item(31, 9).Value = 1538.4517
$ws.Cells.Item(31, 10).Value = 5285.95
$ws.Cells.Item(31, 11).Value = 1538.4517
$ws.Cells.Item(31, 12).Value = 5285.95
$ws.Cells.Item(31, 13).Value = -1243.4517
$ws.Cells.Item(31, 14).Value = -5875.95
$ws.Cells.Item(34, 8).Value = 3008.0588
$ws.Cells.Item(34, 9).Value = 1538.4517
$ws.Cells.Item(34, 10).Value = 5285.95
$ws.Cells.Item(34, 11).Value = 1538.4517
$ws.Cells.Item(34, 12).Value = 5285.95
$ws.Cells.Item(34, 13).Value = -1336.4517
$ws.Cells.Item(34, 14).Value = -5689.95
$ws.Cells.Item(58, 8).Value = 29565.777
$ws.Cells.Item(58, 9).Value = 1771.3334
$ws.Cells.Item(58, 11).Value = 1771.3334
$ws.Cells.Item(58, 13).Value = -1568.3334
$ws.Cells.Item(94, 8).Value = 3151.6667
$ws.Cells.Item(94, 9).Value = 2403.25
$ws.Cells.Item(94, 10).Value = 4007
$ws.Cells.Item(94, 11).Value = 2403.25
$ws.Cells.Item(94, 12).Value = 4007
$ws.Cells.Item(94, 13).Value = -1952.25
$ws.Cells.Item(94, 14).Value = -4909
$ws.Cells.Item(99, 8).Value = 3175.6924
$ws.Cells.Item(99, 9).Value = 2709.8823
$ws.Cells.Item(99, 10).Value = 4055.5557
$ws.Cells.Item(99, 11).Value = 2709.8823
$ws.Cells.Item(99, 12).Value = 4055.5557
$ws.Cells.Item(99, 13).Value = -1211.8823
$ws.Cells.Item(99, 14).Value = -7051.5557
$ws.Cells.Item(126, 8).Value = 3175.6924
$ws.Cells.Item(126, 9).Value = 2709.8823
$ws.Cells.Item(126, 10).Value = 4055.5557
$ws.Cells.Item(126, 11).Value = 8129.646900000001
$ws.Cells.Item(126, 12).Value = 12166.6671
$ws.Cells.Item(126, 13).Value = -5659.646900000001
$ws.Cells.Item(126, 14).Value = -17106.6671
$ws.Cells.Item(132, 8).Value = 1907.6123
$ws.Cells.Item(132, 9).Value = 1503.3684
$ws.Cells.Item(132, 11).Value = 4510.1052
$ws.Cells.Item(132, 13).Value = -1980.1052
$ws.Cells.Item(136, 8).Value = 29565.777
$ws.Cells.Item(136, 9).Value = 1771.3334
$ws.Cells.Item(136, 11).Value = 5314.0002
$ws.Cells.Item(136, 13).Value = -2764.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value = 8149.7144
$ws.Cells.Item(88, 10).Value = 8149.7144
$ws.Cells.Item(88, 12).Value = 24449.1432
$ws.Cells.Item(88, 14).Value = -25305.1432
$ws.Cells.Item(91, 8).Value = 8149.7144
$ws.Cells.Item(91, 10).Value = 8149.7144
$ws.Cells.Item(91, 12).Value = 24449.1432
$ws.Cells.Item(91, 14).Value = -27413.1432
$ws.Cells.Item(131, 8).Value = 694.65656
$ws.Cells.Item(131, 10).Value = 723.0111000000001
$ws.Cells.Item(131, 12).Value = 2169.0333
$ws.Cells.Item(131, 14).Value = -12249.0333
$ws.Cells.Item(132, 8).Value = 595.625
$ws.Cells.Item(132, 9).Value = 593.3333
$ws.Cells.Item(132, 10).Value = 602.5
$ws.Cells.Item(132, 11).Value = 5339.9997
$ws.Cells.Item(132, 12).Value = 5422.5
$ws.Cells.Item(132, 13).Value = -2809.9997
$ws.Cells.Item(132, 14).Value = -10482.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 15410.391
$ws.Cells.Item(132, 9).Value = 3067.0334
$ws.Cells.Item(132, 10).Value = 49074.09
$ws.Cells.Item(132, 11).Value = 9201.100199999999
$ws.Cells.Item(132, 12).Value = 147222.27
$ws.Cells.Item(132, 13).Value = -6671.100199999999
$ws.Cells.Item(132, 14).Value = -152282.27
$ws.Cells.Item(134, 8).Value = 26494.5
$ws.Cells.Item(134, 10).Value = 26494.5
$ws.Cells.Item(134, 12).Value = 79483.5
$ws.Cells.Item(134, 14).Value = -84553.5
$ws.Cells.Item(135, 8).Value = 39773.6
$ws.Cells.Item(135, 10).Value = 39773.6
$ws.Cells.Item(135, 12).Value = 39773.6
$ws.Cells.Item(135, 14).Value = -49913.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4885.0713
$ws.Cells.Item(61, 9).Value = 1743.5555
$ws.Cells.Item(61, 11).Value = 1743.5555
$ws.Cells.Item(61, 13).Value = -1541.5555
$ws.Cells.Item(100, 8).Value = 2281.75
$ws.Cells.Item(100, 9).Value = 1516
$ws.Cells.Item(100, 11).Value = 1516
$ws.Cells.Item(100, 13).Value = -975
$ws.Cells.Item(113, 8).Value = 4885.0713
$ws.Cells.Item(113, 9).Value = 1743.5555
$ws.Cells.Item(113, 11).Value = 1743.5555
$ws.Cells.Item(113, 13).Value = 426.4445000000001
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3287.6667
$ws.Cells.Item(81, 10).Value = 3931.5
$ws.Cells.Item(81, 12).Value = 7863
$ws.Cells.Item(81, 14).Value = -9985
$ws.Cells.Item(84, 8).Value = 3287.6667
$ws.Cells.Item(84, 10).Value = 3931.5
$ws.Cells.Item(84, 12).Value = 39315
$ws.Cells.Item(84, 14).Value = -49923
$ws.Cells.Item(113, 8).Value = 1011.64703
$ws.Cells.Item(113, 9).Value = 1011.64703
$ws.Cells.Item(113, 11).Value = 3034.94109
$ws.Cells.Item(113, 13).Value = -864.9410899999998
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()
